$d = $word.ActiveDocument

# Locate the paragraph that ends with "Tarık Eren Tosun 1210606015" (it
# currently also carries the _GoBack bookmark at its end).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*1210606015*") {
        $target = $p
        break
    }
}

# Insert a brand-new list paragraph right after it; it inherits the
# same list style / level (ListeParagraf, ilvl=2, numId=1) from $target.
$target.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($i + 1)

# Build the new paragraph's content as OOXML so we can split the name
# "Hasgeçkin" into its own run wrapped with spellcheck proofErr markers,
# exactly as Word would when it flags a word unknown to the dictionary.
$cc = [char]0x00e7
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:pPr><w:pStyle w:val="ListeParagraf"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
       '<w:r><w:t xml:space="preserve">Mert </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>Hasge' + $cc + 'kin</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> 1210606072</w:t></w:r>' +
       '<w:bookmarkStart w:id="99" w:name="_GoBack"/><w:bookmarkEnd w:id="99"/>' +
       '</w:p>' +
       '</w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$newPara.Range.InsertXML($xml) | Out-Null

# The _GoBack bookmark now exists twice (the original one still sitting
# in the "Tarık" paragraph, plus the freshly inserted one). Word keeps
# only a single _GoBack marking the most recent edit, so drop the stale
# one left behind in the previous paragraph.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()
